$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions) - update "想去人数" (want-to-go count) column F
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 221
$ws.Range("F3").Value = 2505
$ws.Range("F6").Value = 1320
$ws.Range("F11").Value = 28
$ws.Range("F12").Value = 1790
$ws.Range("F14").Value = 1873
$ws.Range("F16").Value = 1043
$ws.Range("F17").Value = 56
$ws.Range("F19").Value = 1610
$ws.Range("F23").Value = 2407
$ws.Range("F24").Value = 452
$ws.Range("F26").Value = 1038
$ws.Range("F27").Value = 4596
$ws.Range("F31").Value = 173
$ws.Range("F34").Value = 1000

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 30
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 12
$ws.Range("F20").Value = 17
$ws.Range("F21").Value = 17
$ws.Range("F32").Value = 481
$ws.Range("F49").Value = 48

# Sheet 3: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2570
$ws.Range("F3").Value = 2568
$ws.Range("F4").Value = 9599
$ws.Range("F5").Value = 176
$ws.Range("F9").Value = 3114
$ws.Range("F10").Value = 625
$ws.Range("F11").Value = 898
$ws.Range("F12").Value = 322
$ws.Range("F14").Value = 61
$ws.Range("F15").Value = 16
$ws.Range("F16").Value = 315

# Sheet 4: 全部类型 (All Types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2570
$ws.Range("F3").Value = 176
$ws.Range("F4").Value = 2505
$ws.Range("F6").Value = 3114
$ws.Range("F7").Value = 625
$ws.Range("F8").Value = 898
$ws.Range("F12").Value = 61
$ws.Range("F13").Value = 61
$ws.Range("F14").Value = 1320
$ws.Range("F15").Value = 28
$ws.Range("F16").Value = 16
$ws.Range("F17").Value = 1790
$ws.Range("F19").Value = 12
$ws.Range("F20").Value = 1873
$ws.Range("F21").Value = 1043
$ws.Range("F22").Value = 56
$ws.Range("F24").Value = 1610
$ws.Range("F27").Value = 17
$ws.Range("F31").Value = 2407
$ws.Range("F32").Value = 452
$ws.Range("F35").Value = 1038
$ws.Range("F37").Value = 315
$ws.Range("F41").Value = 481
$ws.Range("F46").Value = 173
$ws.Range("F51").Value = 1000
